$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 7.408999999999999
$ws.Range("B7").Value = 5.24
$ws.Range("E7").Value = 16.348
$ws.Range("E15").Value = 16.113
$ws.Range("B16").Value = 5.497999999999999
$ws.Range("E21").Value = 16.557
$ws.Range("E22").Value = 16.559
$ws.Range("E23").Value = 16.527
$ws.Range("B28").Value = 6.034
$ws.Range("B29").Value = 5.427
$ws.Range("B32").Value = 6.453
$ws.Range("E34").Value = 16.805
$ws.Range("B40").Value = 9.327999999999999
$ws.Range("E43").Value = 17.117
$ws.Range("E45").Value = 16.92
$ws.Range("E50").Value = 16.535
$ws.Range("E51").Value = 16.808
$ws.Range("B52").Value = 4.944000000000001
$ws.Range("B57").Value = 5.090999999999999
$ws.Range("B66").Value = 5.013
$ws.Range("E66").Value = 17.447
$ws.Range("E67").Value = 17.43
$ws.Range("E79").Value = 17.07
$ws.Range("E84").Value = 16.608
$ws.Range("E92").Value = 17.885
$ws.Range("E97").Value = 16.872
$ws.Range("B100").Value = 5.586999999999999
